$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value.
# NumberFormat is forced to text ("@") before assignment and the original
# style is restored afterward so numeric-looking strings (e.g. "214.56")
# are stored as literal text, matching the source inline-string cells.
$changes = @(
    ,@('D2', '25.892.07')
    ,@('E2', '  +0.44%  ')
    ,@('D3', '1.633.05')
    ,@('E3', '  +0.23%  ')
    ,@('E4', '  +0.67%  ')
    ,@('D5', '214.56')
    ,@('E5', '  +0.07%  ')
    ,@('E7', '  +0.57%  ')
    ,@('E8', '  -0.24%  ')
    ,@('E9', '  -0.16%  ')
    ,@('D10', '19.52')
    ,@('E10', '  -0.70%  ')
    ,@('E11', '  -0.13%  ')
    ,@('D12', '1.859.93')
    ,@('E12', '  +0.29%  ')
    ,@('B13', 'Polkadot')
    ,@('C13', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot')
    ,@('D13', '4.24')
    ,@('E13', '  -0.30%  ')
    ,@('B14', 'WrappedEther')
    ,@('C14', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth')
    ,@('D14', '1.625.90')
    ,@('E14', '  +0.25%  ')
    ,@('D15', '0.543')
    ,@('E15', '  -2.08%  ')
    ,@('E16', '  -0.10%  ')
    ,@('E17', '  -1.13%  ')
    ,@('D18', '25.909.24')
    ,@('E18', '  +0.55%  ')
    ,@('E19', '  +0.49%  ')
    ,@('D20', '193.01')
    ,@('E20', '  +0.93%  ')
    ,@('D21', '4.38')
    ,@('E21', '  -1.71%  ')
    ,@('D22', '9.94')
    ,@('E22', '  +0.12%  ')
    ,@('E23', '  -0.73%  ')
    ,@('E24', '  -0.69%  ')
    ,@('D25', '143.10')
    ,@('E25', '  +0.82%  ')
    ,@('E26', '  +0.35%  ')
    ,@('E27', '  +2.01%  ')
    ,@('E28', '  +0.05%  ')
    ,@('E29', '  -0.31%  ')
    ,@('E30', '  -0.13%  ')
    ,@('E31', '  +0.57%  ')
    ,@('E32', '  -0.90%  ')
    ,@('E33', '  -0.35%  ')
    ,@('E34', '  -0.57%  ')
    ,@('E35', '  +2.13%  ')
    ,@('E36', '  -0.52%  ')
    ,@('D37', '1.137.98')
    ,@('E37', '  -0.43%  ')
    ,@('D38', '0.549')
    ,@('E38', '  +1.02%  ')
    ,@('E39', '  -0.93%  ')
    ,@('E40', '  +0.26%  ')
    ,@('E41', '  +0.59%  ')
    ,@('D42', '0.804')
    ,@('E42', '  -0.23%  ')
    ,@('D43', '99.15')
    ,@('E43', '  -1.55%  ')
    ,@('E44', '  -2.65%  ')
    ,@('D45', '1.769.44')
    ,@('E45', '  +0.30%  ')
    ,@('B46', 'BabyDogeCoin')
    ,@('C46', 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge')
    ,@('D46', '0.0₆0111')
    ,@('E46', '  -0.47%  ')
    ,@('B47', 'Aave')
    ,@('C47', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave')
    ,@('D47', '56.25')
    ,@('E47', '  +1.71%  ')
    ,@('B48', 'Cronos')
    ,@('C48', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro')
    ,@('D48', '0.0529')
    ,@('E48', '  +3.40%  ')
    ,@('B49', 'RenderToken')
    ,@('C49', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr')
    ,@('D49', '1.46')
    ,@('E49', '  -0.15%  ')
    ,@('B50', 'Mantle')
    ,@('C50', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt')
    ,@('D50', '0.416')
    ,@('E50', '  -0.01%  ')
    ,@('B51', 'EnergySwap')
    ,@('C51', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens')
    ,@('D51', '7.63')
    ,@('E51', '  +1.53%  ')
)

foreach ($change in $changes) {
    $cellRef = $change[0]
    $newValue = $change[1]
    $range = $ws.Range($cellRef)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $newValue
    $range.Style = $origStyle
}

